# Update the NATMI ligand-receptor edge table (Mst1-Mst1r) with newly computed
# TPM-derived values. Ligand/Receptor expression stats (columns G-P) and the
# derived Edge weights/specificities (columns Q-T, which equal G*M, H*N, I*O, J*P
# respectively) all shift together per the updated source TPM data; only the text
# columns (A-D), cell counts (E) and detection-rate numerator bases (F) stay fixed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4312776666666667
$ws.Range("H2").Value = 1.293833
$ws.Range("I2").Value = 0.1323333139342157
$ws.Range("J2").Value = 0.1323333139342157
$ws.Range("M2").Value = 0.1423196666666667
$ws.Range("N2").Value = 0.426959
$ws.Range("O2").Value = 0.03398470032207376
$ws.Range("P2").Value = 0.03398470032207377
$ws.Range("Q2").Value = 0.06137929376077778
$ws.Range("R2").Value = 0.5524136438470001
$ws.Range("S2").Value = 0.00449730801668123
$ws.Range("T2").Value = 0.00449730801668123
$ws.Range("G3").Value = 0.4312776666666667
$ws.Range("H3").Value = 1.293833
$ws.Range("I3").Value = 0.1323333139342157
$ws.Range("J3").Value = 0.1323333139342157
$ws.Range("O3").Value = 0.880663722051367
$ws.Range("P3").Value = 0.880663722051367
$ws.Range("Q3").Value = 1.590554478573444
$ws.Range("R3").Value = 14.314990307161
$ws.Range("S3").Value = 0.1165411488006984
$ws.Range("T3").Value = 0.1165411488006984
$ws.Range("G4").Value = 0.4312776666666667
$ws.Range("H4").Value = 1.293833
$ws.Range("I4").Value = 0.1323333139342157
$ws.Range("J4").Value = 0.1323333139342157
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3139543333333334
$ws.Range("N4").Value = 0.9418630000000001
$ws.Range("O4").Value = 0.0749695680368592
$ws.Range("P4").Value = 0.0749695680368592
$ws.Range("Q4").Value = 0.1354014923198889
$ws.Range("R4").Value = 1.218613430879
$ws.Range("S4").Value = 0.009920971382534232
$ws.Range("T4").Value = 0.009920971382534232
$ws.Range("G5").Value = 0.4312776666666667
$ws.Range("H5").Value = 1.293833
$ws.Range("I5").Value = 0.1323333139342157
$ws.Range("J5").Value = 0.1323333139342157
$ws.Range("M5").Value = 0.04347733333333333
$ws.Range("N5").Value = 0.130432
$ws.Range("O5").Value = 0.01038200958970001
$ws.Range("P5").Value = 0.01038200958970001
$ws.Range("Q5").Value = 0.01875080287288889
$ws.Range("R5").Value = 0.168757225856
$ws.Range("S5").Value = 0.001373885734301809
$ws.Range("T5").Value = 0.001373885734301809
$ws.Range("H6").Value = 4.995672000000001
$ws.Range("I6").Value = 0.5109576205649194
$ws.Range("J6").Value = 0.5109576205649194
$ws.Range("M6").Value = 0.1423196666666667
$ws.Range("N6").Value = 0.426959
$ws.Range("O6").Value = 0.03398470032207376
$ws.Range("P6").Value = 0.03398470032207377
$ws.Range("Q6").Value = 0.2369941246053334
$ws.Range("R6").Value = 2.132947121448001
$ws.Range("S6").Value = 0.01736474161217866
$ws.Range("T6").Value = 0.01736474161217866
$ws.Range("H7").Value = 4.995672000000001
$ws.Range("I7").Value = 0.5109576205649194
$ws.Range("J7").Value = 0.5109576205649194
$ws.Range("O7").Value = 0.880663722051367
$ws.Range("P7").Value = 0.880663722051367
$ws.Range("Q7").Value = 6.141355548269335
$ws.Range("R7").Value = 55.27219993442401
$ws.Range("S7").Value = 0.4499818399372121
$ws.Range("T7").Value = 0.4499818399372121
$ws.Range("H8").Value = 4.995672000000001
$ws.Range("I8").Value = 0.5109576205649194
$ws.Range("J8").Value = 0.5109576205649194
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3139543333333334
$ws.Range("N8").Value = 0.9418630000000001
$ws.Range("O8").Value = 0.0749695680368592
$ws.Range("P8").Value = 0.0749695680368592
$ws.Range("Q8").Value = 0.5228042907706668
$ws.Range("R8").Value = 4.705238616936001
$ws.Range("S8").Value = 0.03830627209889342
$ws.Range("T8").Value = 0.03830627209889342
$ws.Range("H9").Value = 4.995672000000001
$ws.Range("I9").Value = 0.5109576205649194
$ws.Range("J9").Value = 0.5109576205649194
$ws.Range("M9").Value = 0.04347733333333333
$ws.Range("N9").Value = 0.130432
$ws.Range("O9").Value = 0.01038200958970001
$ws.Range("P9").Value = 0.01038200958970001
$ws.Range("Q9").Value = 0.07239949892266667
$ws.Range("R9").Value = 0.6515954903040001
$ws.Range("S9").Value = 0.005304766916635291
$ws.Range("T9").Value = 0.005304766916635291
$ws.Range("G10").Value = 0.7182733333333333
$ws.Range("H10").Value = 2.15482
$ws.Range("I10").Value = 0.2203951139998181
$ws.Range("J10").Value = 0.2203951139998181
$ws.Range("M10").Value = 0.1423196666666667
$ws.Range("N10").Value = 0.426959
$ws.Range("O10").Value = 0.03398470032207376
$ws.Range("P10").Value = 0.03398470032207377
$ws.Range("Q10").Value = 0.1022244213755556
$ws.Range("R10").Value = 0.9200197923800001
$ws.Range("S10").Value = 0.007490061901733103
$ws.Range("T10").Value = 0.007490061901733105
$ws.Range("G11").Value = 0.7182733333333333
$ws.Range("H11").Value = 2.15482
$ws.Range("I11").Value = 0.2203951139998181
$ws.Range("J11").Value = 0.2203951139998181
$ws.Range("O11").Value = 0.880663722051367
$ws.Range("P11").Value = 0.880663722051367
$ws.Range("Q11").Value = 2.648996123548889
$ws.Range("R11").Value = 23.84096511194
$ws.Range("S11").Value = 0.1940939814170152
$ws.Range("T11").Value = 0.1940939814170152
$ws.Range("G12").Value = 0.7182733333333333
$ws.Range("H12").Value = 2.15482
$ws.Range("I12").Value = 0.2203951139998181
$ws.Range("J12").Value = 0.2203951139998181
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.3139543333333334
$ws.Range("N12").Value = 0.9418630000000001
$ws.Range("O12").Value = 0.0749695680368592
$ws.Range("P12").Value = 0.0749695680368592
$ws.Range("Q12").Value = 0.2255050255177778
$ws.Range("R12").Value = 2.02954522966
$ws.Range("S12").Value = 0.01652292649400071
$ws.Range("T12").Value = 0.01652292649400071
$ws.Range("G13").Value = 0.7182733333333333
$ws.Range("H13").Value = 2.15482
$ws.Range("I13").Value = 0.2203951139998181
$ws.Range("J13").Value = 0.2203951139998181
$ws.Range("M13").Value = 0.04347733333333333
$ws.Range("N13").Value = 0.130432
$ws.Range("O13").Value = 0.01038200958970001
$ws.Range("P13").Value = 0.01038200958970001
$ws.Range("Q13").Value = 0.03122860913777778
$ws.Range("R13").Value = 0.28105748224
$ws.Range("S13").Value = 0.002288144187069138
$ws.Range("T13").Value = 0.002288144187069138
$ws.Range("G14").Value = 0.4442506666666666
$ws.Range("H14").Value = 1.332752
$ws.Range("I14").Value = 0.1363139515010468
$ws.Range("J14").Value = 0.1363139515010468
$ws.Range("M14").Value = 0.1423196666666667
$ws.Range("N14").Value = 0.426959
$ws.Range("O14").Value = 0.03398470032207376
$ws.Range("P14").Value = 0.03398470032207377
$ws.Range("Q14").Value = 0.06322560679644444
$ws.Range("R14").Value = 0.569030461168
$ws.Range("S14").Value = 0.004632588791480772
$ws.Range("T14").Value = 0.004632588791480772
$ws.Range("G15").Value = 0.4442506666666666
$ws.Range("H15").Value = 1.332752
$ws.Range("I15").Value = 0.1363139515010468
$ws.Range("J15").Value = 0.1363139515010468
$ws.Range("O15").Value = 0.880663722051367
$ws.Range("P15").Value = 0.880663722051367
$ws.Range("Q15").Value = 1.638398976087111
$ws.Range("R15").Value = 14.745590784784
$ws.Range("S15").Value = 0.1200467518964414
$ws.Range("T15").Value = 0.1200467518964414
$ws.Range("G16").Value = 0.4442506666666666
$ws.Range("H16").Value = 1.332752
$ws.Range("I16").Value = 0.1363139515010468
$ws.Range("J16").Value = 0.1363139515010468
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.3139543333333334
$ws.Range("N16").Value = 0.9418630000000001
$ws.Range("O16").Value = 0.0749695680368592
$ws.Range("P16").Value = 0.0749695680368592
$ws.Range("Q16").Value = 0.1394744218862222
$ws.Range("R16").Value = 1.255269796976
$ws.Range("S16").Value = 0.01021939806143085
$ws.Range("T16").Value = 0.01021939806143085
$ws.Range("G17").Value = 0.4442506666666666
$ws.Range("H17").Value = 1.332752
$ws.Range("I17").Value = 0.1363139515010468
$ws.Range("J17").Value = 0.1363139515010468
$ws.Range("M17").Value = 0.04347733333333333
$ws.Range("N17").Value = 0.130432
$ws.Range("O17").Value = 0.01038200958970001
$ws.Range("P17").Value = 0.01038200958970001
$ws.Range("Q17").Value = 0.01931483431822222
$ws.Range("R17").Value = 0.173833508864
$ws.Range("S17").Value = 0.002288144187069138
$ws.Range("T17").Value = 0.002288144187069138
